# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker rows (16-19) of the "ESTADO DE CUENTA" table are re-sorted so
# that they are grouped by "Periodo Mora" (column E) instead of by worker
# (column C/D). The underlying records are unchanged - only their row
# order (and therefore which row each value lives in) changes:
#
#   before: Martha/1802, Martha/1801, Monica/1802, Monica/1801
#   after : Martha/1801, Monica/1801, Martha/1802, Monica/1802

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the current B16:G19 block (4 rows x 6 cols) into plain variables.
$row16 = @($ws.Range("B16").Value2, $ws.Range("C16").Value2, $ws.Range("D16").Value2, $ws.Range("E16").Value2, $ws.Range("F16").Value2, $ws.Range("G16").Value2)
$row17 = @($ws.Range("B17").Value2, $ws.Range("C17").Value2, $ws.Range("D17").Value2, $ws.Range("E17").Value2, $ws.Range("F17").Value2, $ws.Range("G17").Value2)
$row18 = @($ws.Range("B18").Value2, $ws.Range("C18").Value2, $ws.Range("D18").Value2, $ws.Range("E18").Value2, $ws.Range("F18").Value2, $ws.Range("G18").Value2)
$row19 = @($ws.Range("B19").Value2, $ws.Range("C19").Value2, $ws.Range("D19").Value2, $ws.Range("E19").Value2, $ws.Range("F19").Value2, $ws.Range("G19").Value2)

# New order (equivalent to a stable sort by Periodo Mora ascending):
#   new16 <- old17 (Martha/1801)
#   new17 <- old19 (Monica/1801)
#   new18 <- old16 (Martha/1802)
#   new19 <- old18 (Monica/1802)
$newRow16 = $row17
$newRow17 = $row19
$newRow18 = $row16
$newRow19 = $row18

function Write-RowData($r, $rowData) {
    $ws.Range("B$r").Value = $rowData[0]
    $ws.Range("C$r").Value = $rowData[1]
    $ws.Range("D$r").Value = $rowData[2]
    $ws.Range("E$r").Value = $rowData[3]
    $ws.Range("F$r").Value = $rowData[4]
    $ws.Range("G$r").Value = $rowData[5]
}

Write-RowData 16 $newRow16
Write-RowData 17 $newRow17
Write-RowData 18 $newRow18
Write-RowData 19 $newRow19
